$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.135.61'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +5.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.921.05'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.60%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.21'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5222'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4097'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08524'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.88%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.07'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.11%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.129'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.44'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +9.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.422'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.919.60'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.426'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '95.65'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.89%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06684'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.46'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.135.98'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.36'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.210'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.138.19'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.11'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.13'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.449'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.25'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.084'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1064'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.058'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.640'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02492'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06624'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2209'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.237'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.191'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.916'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6559'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.252'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.64'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.77%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.29'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.766'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.080'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.50%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.72'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.170'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +11.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.85'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.28%  '
